# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> used by the (single) slide master / slides,
#                            originally the "Integral" colour scheme.
#   ppt/theme/theme2.xml -> used only by the notes master, originally the
#                            stock "Office Theme" colour scheme.
#
# The authored change swaps the two themes' colour schemes: the slide
# design now uses the stock "Office Theme" palette while the notes master
# keeps "Integral" (the font scheme / format scheme are byte-identical
# between the two themes, so only the 12 theme colours actually differ).
#
# Apply the new palette to the presentation's theme through the
# SlideMaster's ThemeColorScheme. PowerPoint exposes the theme colours in
# the fixed order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$cs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$cs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$cs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$cs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$cs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$cs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$cs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$cs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$cs.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$cs.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$cs.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
